$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# New user rows to append (id, uin, name, email, mobile)
$newRows = @(
    @{ Id = 110021; Uin = 7316931025; Name = ("Magdalena" + $nbsp + "Weber"); Email = "magdalena.weber@xyz.com"; Mobile = 932122450 },
    @{ Id = 110022; Uin = 9137847236; Name = ("Adrienne" + $nbsp + "Hoffman"); Email = "adrienne.hoffman@xyz.com"; Mobile = 848488000 },
    @{ Id = 110023; Uin = 8428758532; Name = ("Adrienne" + $nbsp + "Mcgee"); Email = "adrienne.mcgee@xyz.com"; Mobile = 894773246 },
    @{ Id = 110024; Uin = 9804209494; Name = ("Amare" + $nbsp + "Coleman"); Email = "amare.coleman@xyz.com"; Mobile = 956554588 },
    @{ Id = 110025; Uin = 7105248214; Name = ("Dawson" + $nbsp + "Ibarra"); Email = "dawson.ibarra@xyz.com"; Mobile = 765455583 },
    @{ Id = 110026; Uin = 9316557128; Name = ("Elvis" + $nbsp + "Mcmillan"); Email = "elvis.mcmillan@xyz.com"; Mobile = 884282274 },
    @{ Id = 110027; Uin = 8103486949; Name = ("Steve" + $nbsp + "George"); Email = "steve.george@xyz.com"; Mobile = 971073663 },
    @{ Id = 110028; Uin = 9601932866; Name = ("Colton" + $nbsp + "Elliott"); Email = "colton.elliott@xyz.com"; Mobile = 809908673 },
    @{ Id = 110029; Uin = 9317596765; Name = ("Carolyn" + $nbsp + "Rodriguez"); Email = "carolyn.rodriguez@xyz.com"; Mobile = 818876429 }
)

$startRow = 22

# Populate column-by-column (A, B, C, D, E, F, G, H, I, J, K) so that new
# shared-string entries are interned in the same order Excel produced them
# in (all names, then all emails, etc.), matching the original authoring
# pattern of this sheet.

# Column A - id
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.Id
    $r++
}

# Column B - uin
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 2).Value = $row.Uin
    $r++
}

# Column C - name
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 3).Value = $row.Name
    $r++
}

# Column D - email
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 4).Value = $row.Email
    $r++
}

# Column E - mobile
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 5).Value = $row.Mobile
    $r++
}

# Column F - status_code
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 6).Value = "ACT"
    $r++
}

# Column G - lang_code
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 7).Value = "eng"
    $r++
}

# Column H - last_login_method
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 8).Value = "PWD"
    $r++
}

# Column I - is_active
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 9).Value = $true
    $r++
}

# Column J - cr_by
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 10).Value = "superadmin"
    $r++
}

# Column K - cr_dtimes
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 11).Value = "now()"
    $r++
}

# Match formatting of the preceding data rows: is_active column (I) is
# left-aligned, matching every other data row in the table.
$ws.Range("I22:I30").HorizontalAlignment = -4131  # xlLeft

# Update view: scroll so row 16 is the top row, select the newly added block.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A22:K30").Select()
